$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=15.30803754805065; "C"=7.466838059971029; "E"=11.63119828460207; "F"=16.86991607391233; "G"=43.34839462607674; "H"=17.96528885461717; "K"=11.62206511110484; "L"=9.879110084296308; "N"=21.05370929473121 }
    3 = @{ "B"=15.11697530907573; "C"=7.436619688945807; "E"=11.63632060300179; "F"=15.89584955866808; "G"=43.40566878627979; "H"=18.01548792079482; "K"=11.49461062204216; "L"=9.8690195119026; "N"=21.12182643385264 }
    4 = @{ "B"=15.00214567654348; "C"=7.417666176567137; "E"=11.64148435617159; "F"=15.26997757108491; "G"=43.45257127351685; "H"=18.04923104990212; "K"=11.41825837136148; "L"=9.864581121215602; "N"=21.16560678525283 }
    5 = @{ "B"=14.9560323016867; "C"=7.409843296101898; "E"=11.64409671408608; "F"=15.00819731993403; "G"=43.47462602694754; "H"=18.06371534946115; "K"=11.38765787214931; "L"=9.863216080908236; "N"=21.18394089390627 }
    6 = @{ "B"=14.94841793575135; "C"=7.408538359250096; "E"=11.64456119213344; "F"=14.96433081551593; "G"=43.47846556511985; "H"=18.06616474828299; "K"=11.38260869667855; "L"=9.863016258072962; "N"=21.18701509068231 }
    7 = @{ "B"=15.00152094633229; "C"=7.417561074359695; "E"=11.64151752960194; "F"=15.26647399323137; "G"=43.45285681474633; "H"=18.04942342051258; "K"=11.4178435576249; "L"=9.864560913390989; "N"=21.16585204646296 }
    8 = @{ "B"=15.24167793096979; "C"=7.456502382802962; "E"=11.63254583322109; "F"=16.53996406344768; "G"=43.3657016479451; "H"=17.98199086682603; "K"=11.5777443730657; "L"=9.875267268543395; "N"=21.0767909132015 }
    9 = @{ "B"=15.72981084228364; "C"=7.529645138443356; "E"=11.63094212358256; "F"=19.00274580682531; "G"=43.28827821481014; "H"=17.8729627275854; "K"=11.90486738284968; "L"=9.910117787858908; "N"=20.91759978085068 }
    10 = @{ "B"=16.0955079849236; "C"=7.581351514942878; "E"=11.63946654085591; "F"=20.67494806633232; "G"=43.28880785969233; "H"=17.80705117165894; "K"=12.15134036682225; "L"=9.944036553951314; "N"=20.80997444013686 }
    11 = @{ "B"=16.26268399296289; "C"=7.604418565906151; "E"=11.64543802632986; "F"=21.3917225636224; "G"=43.30157475773917; "H"=17.78015688643381; "K"=12.26434229587207; "L"=9.961238689186617; "N"=20.76301994336216 }
    12 = @{ "B"=16.32604733792536; "C"=7.613086821564266; "E"=11.64799885943983; "F"=21.65686569030329; "G"=43.30821262127421; "H"=17.77041767406114; "K"=12.30722164203628; "L"=9.968004234393698; "N"=20.7455263682602 }
    13 = @{ "B"=16.31239930099805; "C"=7.61122294778965; "E"=11.64743404124683; "F"=21.60004134736742; "G"=43.30670281957009; "H"=17.77249537996627; "K"=12.29798350455225; "L"=9.9665360257017; "N"=20.74928117418013 }
    14 = @{ "B"=16.26789612371632; "C"=7.605133053945718; "E"=11.6456427132931; "F"=21.4136618050453; "G"=43.30208471140916; "H"=17.779346710978; "K"=12.26786846183136; "L"=9.961790275160382; "N"=20.76157498993911 }
    15 = @{ "B"=16.24064235156908; "C"=7.601394085979398; "E"=11.64458443680281; "F"=21.29868154950795; "G"=43.29949086625616; "H"=17.78360133483303; "K"=12.24943244727367; "L"=9.9589160097388; "N"=20.76914265670387 }
    16 = @{ "B"=16.08459393096965; "C"=7.579834765731083; "E"=11.63911829608028; "F"=20.62722412089977; "G"=43.28822578378717; "H"=17.80887101218347; "K"=12.14396982635354; "L"=9.942947733978231; "N"=20.8130832646196 }
    17 = @{ "B"=15.9890296798127; "C"=7.566491440082357; "E"=11.63630022787213; "F"=20.20408069597325; "G"=43.28452522863857; "H"=17.82516497251105; "K"=12.07946975618277; "L"=9.93360332088314; "N"=20.84055201522182 }
    18 = @{ "B"=15.93414424087428; "C"=7.558774021122412; "E"=11.63487651982336; "F"=19.95656407809801; "G"=43.28357574666638; "H"=17.83482755000402; "K"=12.0424561528062; "L"=9.928395700207213; "N"=20.85654008802886 }
    19 = @{ "B"=15.91557667138143; "C"=7.556153745767998; "E"=11.63442838366586; "F"=19.87204792380568; "G"=43.28345666874158; "H"=17.83814903839065; "K"=12.02993978697205; "L"=9.926661277261987; "N"=20.86198583417064 }
    20 = @{ "B"=15.99919475705481; "C"=7.5679162883017; "E"=11.63657981997184; "F"=20.24955283636154; "G"=43.28479711866266; "H"=17.82340035769582; "K"=12.08632736900178; "L"=9.934580785448457; "N"=20.83760839125419 }
    21 = @{ "B"=16.28096670174324; "C"=7.606923625634203; "E"=11.646160752722; "F"=21.46857628470577; "G"=43.30339221160094; "H"=17.77732222332967; "K"=12.27671189764222; "L"=9.963177420601268; "N"=20.75795621676243 }
    22 = @{ "B"=16.4654285996427; "C"=7.632027545589247; "E"=11.65416768014923; "F"=22.22866616901552; "G"=43.32605610351944; "H"=17.74980207651408; "K"=12.40163441870914; "L"=9.983330980018383; "N"=20.70757164757187 }
    23 = @{ "B"=16.36696925109074; "C"=7.618665220104375; "E"=11.6497350920527; "F"=21.82633154458858; "G"=43.312997937377; "H"=17.76425242148914; "K"=12.3349282213805; "L"=9.972441875634736; "N"=20.73431018344535 }
    24 = @{ "B"=15.99459894743666; "C"=7.567272258128164; "E"=11.63645280439962; "F"=20.22900810905287; "G"=43.28467052783294; "H"=17.82419722128848; "K"=12.08322682639088; "L"=9.934138360734874; "N"=20.83893859323264 }
    25 = @{ "B"=15.59626981854826; "C"=7.510209920728368; "E"=11.62966807067825; "F"=18.34778573295695; "G"=43.29916560218716; "H"=17.89996865047533; "K"=11.81513901697719; "L"=9.899219086268287; "N"=20.9590198170167 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
